# Update coefficient table values with re-run regression results ("more adjustments").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.1016322605355135
$ws.Range("B3").Value = 0.1331043532377512
$ws.Range("H3").Value = 0.2347366137732647
$ws.Range("B4").Value = 0.1580679086825927
$ws.Range("H4").Value = 0.2597001692181062
$ws.Range("B5").Value = 0.04376156664135121
$ws.Range("C5").Value = 0.002529839505213326
$ws.Range("D5").Value = 8.761641952528331
$ws.Range("E5").Value = 0.01542800872070308
$ws.Range("F5").Value = 0.03879675392540394
$ws.Range("G5").Value = 0.04872637935729793
$ws.Range("H5").Value = 0.1453938271768647
$ws.Range("B6").Value = 0.02783253679909668
$ws.Range("C6").Value = 0.001766389746154192
$ws.Range("D6").Value = 4.776274590672215
$ws.Range("E6").Value = 0.01275749563726479
$ws.Range("F6").Value = 0.02436808887859849
$ws.Range("G6").Value = 0.03129698471959536
$ws.Range("H6").Value = 0.1294647973346101
$ws.Range("B7").Value = 0.01834357672657711
$ws.Range("C7").Value = 0.001438133890871352
$ws.Range("D7").Value = 3.184669944891964
$ws.Range("E7").Value = 0.01534245702020938
$ws.Range("F7").Value = 0.01552136312821852
$ws.Range("G7").Value = 0.02116579032493603
$ws.Range("H7").Value = 0.1199758372620906
$ws.Range("B8").Value = 0.01162570027688834
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").Value = 0.1132579608124018
$ws.Range("B9").Value = 0.01218648742028626
$ws.Range("C9").Value = 0.00105278393511706
$ws.Range("D9").Value = 1.383215647023965
$ws.Range("E9").Value = 0.006056334445960623
$ws.Range("F9").Value = 0.01012114568143264
$ws.Range("G9").Value = 0.01425182915913969
$ws.Range("H9").Value = 0.1138187479557997
$ws.Range("B10").Value = 0.01164651242353386
$ws.Range("C10").Value = 0.001178678676683818
$ws.Range("D10").Value = 1.402255000863588
$ws.Range("E10").Value = 0.00763873425262017
$ws.Range("F10").Value = 0.009334718847606192
$ws.Range("G10").Value = 0.01395830599946176
$ws.Range("H10").Value = 0.1132787729590473
$ws.Range("B11").Value = 0.03383545647285402
$ws.Range("H11").Value = 0.1354677170083675
$ws.Range("B12").Value = 0.05926678375606758
$ws.Range("H12").Value = 0.1608990442915811
$ws.Range("B13").Value = 0.07702324577357501
$ws.Range("H13").Value = 0.1786555063090885
$ws.Range("B14").Value = 0.08557982281271152
$ws.Range("H14").Value = 0.187212083348225
$ws.Range("B15").Value = 0.09137061353124869
$ws.Range("H15").Value = 0.1930028740667621
$ws.Range("B16").Value = 0.09614514120463856
$ws.Range("H16").Value = 0.197777401740152
$ws.Range("B17").Value = 0.0966519284594832
$ws.Range("H17").Value = 0.1982841889949967
$ws.Range("B18").Value = -0.1016322605355135
$ws.Range("C18").Value = 0.008184724085594378
$ws.Range("D18").Value = -19.78194440735432
$ws.Range("E18").Value = 0.0241487565814973
$ws.Range("F18").Value = -0.1177433714158163
$ws.Range("G18").Value = -0.08552114965521068
$ws.Range("B19").Value = 0.09991613980274094
$ws.Range("H19").Value = 0.2015484003382544
$ws.Range("B20").Value = 0.1030089044794883
$ws.Range("H20").Value = 0.2046411650150018
$ws.Range("B21").Value = 0.1069993869585494
$ws.Range("C21").Value = 0.006737559737895482
$ws.Range("D21").Value = 27.23521515254481
$ws.Range("E21").Value = 0.04266115757544881
$ws.Range("F21").Value = 0.09374248010681115
$ws.Range("G21").Value = 0.1202562938102871
$ws.Range("H21").Value = 0.2086316474940628
$ws.Range("B22").Value = 0.1112317171313762
$ws.Range("H22").Value = 0.2128639776668896
$ws.Range("B23").Value = 0.11472946700891
$ws.Range("C23").Value = 0.006933558217167623
$ws.Range("D23").Value = 28.76314926053296
$ws.Range("E23").Value = 0.03970665460976443
$ws.Range("F23").Value = 0.1011106768971863
$ws.Range("G23").Value = 0.1283482571206334
$ws.Range("H23").Value = 0.2163617275444235
$ws.Range("B24").Value = 0.1191865821360233
$ws.Range("C24").Value = 0.006561567039717055
$ws.Range("D24").Value = 30.11727319675678
$ws.Range("E24").Value = 0.03672389257228235
$ws.Range("F24").Value = 0.1062944953281114
$ws.Range("G24").Value = 0.1320786689439353
$ws.Range("H24").Value = 0.2208188426715368
$ws.Range("B25").Value = 0.1242988019061651
$ws.Range("C25").Value = 0.006690923673680792
$ws.Range("D25").Value = 30.792460216712
$ws.Range("E25").Value = 0.03913274174929359
$ws.Range("F25").Value = 0.1111473941228753
$ws.Range("G25").Value = 0.1374502096894548
$ws.Range("H25").Value = 0.2259310624416785
$ws.Range("B26").Value = 0.12590755635499
$ws.Range("C26").Value = 0.006783264532057268
$ws.Range("D26").Value = 31.40931333238305
$ws.Range("E26").Value = 0.03079473828585882
$ws.Range("F26").Value = 0.1125789749865675
$ws.Range("G26").Value = 0.1392361377234135
$ws.Range("H26").Value = 0.2275398168905035
$ws.Range("B27").Value = 0.1313646774139538
$ws.Range("C27").Value = 0.00646597891310742
$ws.Range("D27").Value = 31.78002835819297
$ws.Range("E27").Value = 0.03787074213135282
$ws.Range("F27").Value = 0.1186646448267631
$ws.Range("G27").Value = 0.1440647100011441
$ws.Range("H27").Value = 0.2329969379494672
$ws.Range("B28").Value = 0.1385421202412916
$ws.Range("C28").Value = 0.006929313523036389
$ws.Range("D28").Value = 30.33036311964915
$ws.Range("E28").Value = 0.06298901182335746
$ws.Range("F28").Value = 0.1249367251303649
$ws.Range("G28").Value = 0.1521475153522188
$ws.Range("H28").Value = 0.2401743807768051
$ws.Range("B29").Value = 0.0148064697366559
$ws.Range("C29").Value = 0.001128555487046371
$ws.Range("D29").Value = 2.379300267459233
$ws.Range("E29").Value = 0.02966488404541591
$ws.Range("F29").Value = 0.01256775377464908
$ws.Range("G29").Value = 0.01704518569866261
$ws.Range("H29").Value = 0.1164387302721694